$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '67.765.11'
$ws.Range('E2').Value = '  -2.18%  '
$ws.Range('D3').Value = '2.677.43'
$ws.Range('E3').Value = '  -2.72%  '
$ws.Range('E4').Value = '  +0.15%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '601.08'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.67%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '167.88'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.59%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.546'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -0.25%  '
$ws.Range('D9').Value = '2.676.84'
$ws.Range('E9').Value = '  -2.76%  '
$ws.Range('E10').Value = '  +1.14%  '
$ws.Range('E11').Value = '  +1.18%  '
$ws.Range('E12').Value = '  -0.40%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '5.22'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -2.16%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '27.92'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -3.11%  '
$ws.Range('D15').Value = '3.171.75'
$ws.Range('E15').Value = '  -2.50%  '
$ws.Range('E16').Value = '  -3.21%  '
$ws.Range('D17').Value = '67.710.20'
$ws.Range('D18').Value = '2.680.15'
$ws.Range('E18').Value = '  -3.01%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '11.74'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -1.99%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '7.86'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +1.43%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '365.98'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -0.84%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '4.40'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -4.12%  '
$ws.Range('E23').Value = '  -3.36%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.03'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -5.28%  '
$ws.Range('E25').Value = '  +0.06%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '71.03'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -4.43%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '10.17'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +1.01%  '
$ws.Range('D28').Value = '2.833.07'
$ws.Range('E28').Value = '  -1.53%  '
$ws.Range('E29').Value = '  -4.54%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.999'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -0.07%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '553.37'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -8.40%  '
$ws.Range('E32').Value = '  -3.87%  '
$ws.Range('E34').Value = '  -1.83%  '
$ws.Range('E35').Value = '  -2.03%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.999'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -0.09%  '
$ws.Range('E37').Value = '  -5.37%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '19.52'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -3.37%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '155.31'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -4.90%  '
$ws.Range('E40').Value = '  -2.82%  '
$ws.Range('B41').Value = 'Stacks'
$ws.Range('C41').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.84'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -4.64%  '
$ws.Range('B42').Value = 'RenderToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '5.31'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -4.21%  '
$ws.Range('E43').Value = '  -0.62%  '
$ws.Range('E44').Value = '  -7.67%  '
$ws.Range('E45').Value = '  +0.03%  '
$ws.Range('E46').Value = '  -0.99%  '
$ws.Range('E47').Value = '  -6.44%  '
$ws.Range('E48').Value = '  -3.68%  '
$ws.Range('E49').Value = '  -3.63%  '
$ws.Range('E50').Value = '  -2.64%  '
$ws.Range('E51').Value = '  -4.32%  '
